$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the "Play Champions of Rome Free Slot: Unique Bonus
#    Features" heading.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Right before the closing "Prompt: ..." paragraph, insert a new
#    paragraph containing the bold heading text
#    "Play Champions of Rome Free Slot: Unique Bonus Features".
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($lastIndex)

$insertionPoint = $promptPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore() | Out-Null

$newHeadingIndex = $lastIndex
$newHeadingRange = $d.Paragraphs.Item($newHeadingIndex).Range
$newHeadingRange.Collapse(1)

$headingXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Champions of Rome Free Slot: Unique Bonus Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newHeadingRange.InsertXML($headingXml) | Out-Null

# ------------------------------------------------------------------
# 3) Replace the old "Prompt: ..." text (now the last paragraph) with
#    the meta-description body text, keeping its existing (italic)
#    run formatting untouched.
# ------------------------------------------------------------------
$finalIndex = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalIndex)
$finalPara.Range.Find.Execute(
    "Prompt: Please create a cartoon-style feature image for the game ""Champions of Rome"" featuring a happy Maya warrior with glasses. The image should be eye-catching and engaging, depicting the combination of ancient Roman and Maya cultures in a fun and exciting way.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Champions of Rome is an exciting gladiator game with a user-friendly interface and unique bonus features. Play for free and aim for a winning potential of up to 3000x.",
    2) | Out-Null
